$d = $word.ActiveDocument
$d.Content.Find.Execute("go up to 80,000). ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "go up to 48,000). ", 2)
